# Generate Report for Archive
#
# Updates the status text from "Ready for handoff" to "In Translation"
# everywhere it is used (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all
# share the same underlying string), and re-fits the affected status
# columns so their width reflects the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell currently showing the old status text.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The status columns were auto-fit to the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
